$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 569.8
$ws.Range("I2").Value = 590
$ws.Range("K2").Value = 590
$ws.Range("M2").Value = -477

$ws.Range("H32").Value = 787.7778
$ws.Range("I32").Value = 729.4
$ws.Range("K32").Value = 729.4
$ws.Range("M32").Value = -403.4

$ws.Range("H53").Value = 227.94444
$ws.Range("I53").Value = 229.72728
$ws.Range("J53").Value = 225.14285
$ws.Range("K53").Value = 229.72728
$ws.Range("L53").Value = 225.14285
$ws.Range("M53").Value = 407.27272
$ws.Range("N53").Value = -1499.14285

$ws.Range("H55").Value = 491.05554
$ws.Range("I55").Value = 190.5
$ws.Range("J55").Value = 731.5
$ws.Range("K55").Value = 190.5
$ws.Range("L55").Value = 731.5
$ws.Range("M55").Value = 23.5
$ws.Range("N55").Value = -1159.5

$ws.Range("H98").Value = 2442.5715
$ws.Range("I98").Value = 2246.3845
$ws.Range("J98").Value = 4993
$ws.Range("K98").Value = 2246.3845
$ws.Range("L98").Value = 4993
$ws.Range("M98").Value = -748.3845000000001
$ws.Range("N98").Value = -7989

$ws.Range("H113").Value = 3150
$ws.Range("I113").Value = 3060
$ws.Range("J113").Value = 3600
$ws.Range("K113").Value = 3060
$ws.Range("L113").Value = 3600
$ws.Range("M113").Value = 194
$ws.Range("N113").Value = -10108

$ws.Range("H115").Value = 3735.75
$ws.Range("I115").Value = 3735.75
$ws.Range("K115").Value = 11207.25
$ws.Range("M115").Value = -9640.25

$ws.Range("H122").Value = 2442.5715
$ws.Range("I122").Value = 2246.3845
$ws.Range("J122").Value = 4993
$ws.Range("K122").Value = 6739.1535
$ws.Range("L122").Value = 14979
$ws.Range("M122").Value = -4289.1535
$ws.Range("N122").Value = -19879

$ws.Range("H132").Value = 1556.3226
$ws.Range("I132").Value = 1529.0769
$ws.Range("K132").Value = 4587.2307
$ws.Range("M132").Value = -2057.2307

$ws.Range("H137").Value = 3011.7778
$ws.Range("I137").Value = 1964.0769
$ws.Range("J137").Value = 3984.6428
$ws.Range("K137").Value = 5892.2307
$ws.Range("L137").Value = 11953.9284
$ws.Range("M137").Value = -3342.2307
$ws.Range("N137").Value = -17053.9284

$ws.Range("H138").Value = 13134.077
$ws.Range("J138").Value = 13995.25
$ws.Range("L138").Value = 41985.75
$ws.Range("N138").Value = -52265.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14923.585
$ws.Range("I32").Value = 6763.697
$ws.Range("K32").Value = 6763.697
$ws.Range("M32").Value = -6476.697

$ws.Range("H41").Value = 10062
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 10062
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 10062
$ws.Range("N41").Value = -10890
$ws.Range("M41").ClearContents()

$ws.Range("H61").Value = 1263.238
$ws.Range("I61").Value = 1212
$ws.Range("K61").Value = 1212
$ws.Range("M61").Value = -1000

$ws.Range("H122").Value = 402968.16
$ws.Range("I122").Value = 627089.4
$ws.Range("K122").Value = 1881268.2
$ws.Range("M122").Value = -1878818.2

$ws.Range("H132").Value = 1773.415
$ws.Range("I132").Value = 1730.5962
$ws.Range("K132").Value = 5191.7886
$ws.Range("M132").Value = -2661.7886

$ws.Range("H136").Value = 1263.238
$ws.Range("I136").Value = 1212
$ws.Range("K136").Value = 3636
$ws.Range("M136").Value = -1086

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 833.1667
$ws.Range("I94").Value = 619.8
$ws.Range("K94").Value = 619.8
$ws.Range("M94").Value = -168.8

$ws.Range("H105").Value = 4264.0938
$ws.Range("I105").Value = 3747.3333
$ws.Range("K105").Value = 3747.3333
$ws.Range("M105").Value = -2000.3333

$ws.Range("H134").Value = 1321.7693
$ws.Range("J134").Value = 4666
$ws.Range("L134").Value = 13998
$ws.Range("N134").Value = -19068

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H125").Value = 88499.5
$ws.Range("J125").Value = 88499.5
$ws.Range("L125").Value = 88499.5
$ws.Range("N125").Value = -93419.5

$ws.Range("H132").Value = 3325.889
$ws.Range("I132").Value = 2856.682
$ws.Range("K132").Value = 8570.045999999998
$ws.Range("M132").Value = -6040.045999999998

$ws.Range("H134").Value = 3765.4092
$ws.Range("I134").Value = 2420.25
$ws.Range("K134").Value = 7260.75
$ws.Range("M134").Value = -4725.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 33.384617
$ws.Range("J12").Value = 43
$ws.Range("L12").Value = 129
$ws.Range("N12").Value = -475

$ws.Range("H129").Value = 5499
$ws.Range("J129").Value = 3999
$ws.Range("L129").Value = 11997
$ws.Range("N129").Value = -21997

$ws.Range("H131").Value = 1574.6875
$ws.Range("I131").Value = 631
$ws.Range("J131").Value = 1889.25
$ws.Range("K131").Value = 1893
$ws.Range("L131").Value = 5667.75
$ws.Range("M131").Value = 3147
$ws.Range("N131").Value = -15747.75

$ws.Range("H136").Value = 15933.333
$ws.Range("I136").Value = 11400
$ws.Range("J136").Value = 25000
$ws.Range("K136").Value = 34200
$ws.Range("L136").Value = 75000
$ws.Range("M136").Value = -29100
$ws.Range("N136").Value = -85200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8219.909
$ws.Range("J70").Value = 8718
$ws.Range("L70").Value = 8718
$ws.Range("N70").Value = -9258

$ws.Range("H73").Value = 8219.909
$ws.Range("J73").Value = 8718
$ws.Range("L73").Value = 8718
$ws.Range("N73").Value = -10590

$ws.Range("H102").Value = 1587.3704
$ws.Range("I102").Value = 824.5294
$ws.Range("J102").Value = 2884.2
$ws.Range("K102").Value = 824.5294
$ws.Range("L102").Value = 2884.2
$ws.Range("M102").Value = 797.4706
$ws.Range("N102").Value = -6128.2

$ws.Range("H122").Value = 502240.53
$ws.Range("I122").Value = 73157.42999999999
$ws.Range("J122").Value = 1253136
$ws.Range("K122").Value = 219472.29
$ws.Range("L122").Value = 3759408
$ws.Range("M122").Value = -217022.29
$ws.Range("N122").Value = -3764308

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3780.182
$ws.Range("I7").Value = 3199
$ws.Range("J7").Value = 5330
$ws.Range("K7").Value = 3199
$ws.Range("L7").Value = 5330
$ws.Range("M7").Value = -3087
$ws.Range("N7").Value = -5554

$ws.Range("H82").Value = 1904.0588
$ws.Range("I82").Value = 1797.5
$ws.Range("K82").Value = 1797.5
$ws.Range("M82").Value = -1436.5

$ws.Range("H85").Value = 1904.0588
$ws.Range("I85").Value = 1797.5
$ws.Range("K85").Value = 1797.5
$ws.Range("M85").Value = -549.5

$ws.Range("H126").Value = 3780.182
$ws.Range("I126").Value = 3199
$ws.Range("J126").Value = 5330
$ws.Range("K126").Value = 9597
$ws.Range("L126").Value = 15990
$ws.Range("M126").Value = -7127
$ws.Range("N126").Value = -20930

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1660
$ws.Range("I107").Value = 880
$ws.Range("K107").Value = 2640
$ws.Range("M107").Value = -720

$ws.Range("H132").Value = 2774
$ws.Range("I132").Value = 2365.3333
$ws.Range("K132").Value = 7095.999899999999
$ws.Range("M132").Value = -4565.999899999999

$ws.Range("H136").Value = 34918.645
$ws.Range("I136").Value = 2259.739
$ws.Range("J136").Value = 128813
$ws.Range("K136").Value = 6779.217000000001
$ws.Range("L136").Value = 386439
$ws.Range("M136").Value = -4229.217000000001
$ws.Range("N136").Value = -391539
